# Add new columns I (I0) and J (IF) to the prelander sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - use same style as existing header cell H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-23
$values = @(
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 8),
    @(5, 5),
    @(6, 7),
    @(7, 8),
    @(6, 7),
    @(6, 6),
    @(7, 8),
    @(8, 8),
    @(7, 8),
    @(7, 7),
    @(8, 8),
    @(6, 7),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(6, 6),
    @(8, 9),
    @(3, 3),
    @(3, 3)
)

$row = 2
foreach ($pair in $values) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
